$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 644; this shifts existing rows 644-685 down to 645-686,
# matching the target dimension growing from A1:D685 to A1:D686.
$ws.Rows(644).Insert()

# Populate the newly inserted row with the new record.
# Force column A to plain text so "2026/01/15" is stored as a literal string
# (matching the surrounding inlineStr date cells) instead of being
# auto-converted into a date serial number.
$ws.Range("A644").NumberFormat = "@"
$ws.Range("A644").Value = "2026/01/15"
$ws.Range("B644").Value = "木"
$ws.Range("C644").Value = 17
$ws.Range("D644").Value = 201
